# "Generate Report for Handback"
# Fill in the actual handoff/handback completion timestamps for the
# second tracked file (a9227bf9-...) on both the zh-cn and de-de
# report sheets. Row 2 (a136be05-...) already holds its timestamps and
# is left untouched.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-20 04:38:12"
$wsZhCn.Range("H3").Value = "2016-03-20 04:38:51"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-20 04:38:15"
$wsDeDe.Range("H3").Value = "2016-03-20 04:38:56"
